$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

# Each table cell holds one "A op B = C" equation; the table is 20 rows x 5
# columns, read left-to-right/top-to-bottom, matching the order below.
$oldValues = @(
    "24+26=50",
    "42+30=72",
    "20+26=46",
    "16+73=89",
    "69-46=23",
    "27+31=58",
    "81-72=9",
    "79-24=55",
    "41-38=3",
    "48+8=56",
    "67-50=17",
    "63-10=53",
    "77-50=27",
    "91-8=83",
    "14+71=85",
    "96-38=58",
    "79+13=92",
    "34+7=41",
    "9+63=72",
    "24-11=13",
    "13+66=79",
    "5+56=61",
    "80-54=26",
    "21+53=74",
    "0+35=35",
    "80-54=26",
    "64-24=40",
    "98-24=74",
    "18+26=44",
    "70-60=10",
    "35+63=98",
    "43+17=60",
    "96-69=27",
    "58-16=42",
    "39+53=92",
    "67+10=77",
    "20+8=28",
    "66+19=85",
    "11+33=44",
    "5+73=78",
    "83-58=25",
    "96-31=65",
    "46-6=40",
    "58+29=87",
    "29-21=8",
    "10+32=42",
    "9+48=57",
    "98-85=13",
    "6+66=72",
    "91-2=89",
    "70-62=8",
    "43-20=23",
    "49-32=17",
    "43-26=17",
    "8+30=38",
    "5+20=25",
    "71+12=83",
    "66+9=75",
    "78-17=61",
    "28+5=33",
    "26+49=75",
    "60+37=97",
    "36-24=12",
    "58+20=78",
    "12+13=25",
    "55-35=20",
    "46+7=53",
    "79-66=13",
    "27-3=24",
    "65+1=66",
    "37+57=94",
    "22+3=25",
    "94-85=9",
    "13+7=20",
    "89-34=55",
    "75-27=48",
    "47+39=86",
    "27+1=28",
    "61-45=16",
    "90-61=29",
    "36-11=25",
    "92-53=39",
    "89-16=73",
    "93-83=10",
    "64-49=15",
    "4+14=18",
    "26+41=67",
    "58-44=14",
    "40-1=39",
    "74+11=85",
    "78-40=38",
    "21+59=80",
    "33+28=61",
    "20+18=38",
    "89-30=59",
    "33-18=15",
    "7+84=91",
    "10+32=42",
    "71-34=37",
    "21+77=98"
)

$newValues = @(
    "20+38=58",
    "61-48=13",
    "91-81=10",
    "37+14=51",
    "73-0=73",
    "7+60=67",
    "8+51=59",
    "35+42=77",
    "30-19=11",
    "59+36=95",
    "21-18=3",
    "98-42=56",
    "44+24=68",
    "99-66=33",
    "15+20=35",
    "64-52=12",
    "31-8=23",
    "98-8=90",
    "69-68=1",
    "73-64=9",
    "56-8=48",
    "82-4=78",
    "69-57=12",
    "2-0=2",
    "82-41=41",
    "2+56=58",
    "2+69=71",
    "53+42=95",
    "31+9=40",
    "67-41=26",
    "65-0=65",
    "18+51=69",
    "50+24=74",
    "72-30=42",
    "76-36=40",
    "45-36=9",
    "58-21=37",
    "11-3=8",
    "78-30=48",
    "23+17=40",
    "81-53=28",
    "18+36=54",
    "71-27=44",
    "18+38=56",
    "56-39=17",
    "65-19=46",
    "45+15=60",
    "32-1=31",
    "62-9=53",
    "30+54=84",
    "6+15=21",
    "74+2=76",
    "64-41=23",
    "68-49=19",
    "2+64=66",
    "42-3=39",
    "31-1=30",
    "93-55=38",
    "8+61=69",
    "66-38=28",
    "30+5=35",
    "30+58=88",
    "61+10=71",
    "49+11=60",
    "15-4=11",
    "11+15=26",
    "76-2=74",
    "46+47=93",
    "7+82=89",
    "56+6=62",
    "90-20=70",
    "36-20=16",
    "39+50=89",
    "13+23=36",
    "45+22=67",
    "51+29=80",
    "62-20=42",
    "76-33=43",
    "24+54=78",
    "48+40=88",
    "7+90=97",
    "45-41=4",
    "73-53=20",
    "93-31=62",
    "95-76=19",
    "63-41=22",
    "31-18=13",
    "1+57=58",
    "92-44=48",
    "13-5=8",
    "86-41=45",
    "50+45=95",
    "1+57=58",
    "40+5=45",
    "87+5=92",
    "86-59=27",
    "8+87=95",
    "99-14=85",
    "41-31=10",
    "16+46=62"
)

$cols = 5
$updated = 0
for ($i = 0; $i -lt $newValues.Length; $i++) {
    $row = [int]([math]::Floor($i / $cols)) + 1
    $col = ($i % $cols) + 1
    $cell = $t.Cell($row, $col)
    $current = $cell.Range.Text.TrimEnd([char]13, [char]7)
    if ($current -ne $oldValues[$i]) {
        throw "Cell ($row,$col) expected `"$($oldValues[$i])`" but found `"$current`""
    }
    $cell.Range.Text = $newValues[$i]
    $updated++
}
Write-Output "Updated $updated cells"